$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The four data rows (2-5) rotate their weekly data up by one row,
# wrapping the first row's data down to the last row.
# Row 2 <- old Row 3, Row 3 <- old Row 4, Row 4 <- old Row 5, Row 5 <- old Row 2.

$rows = @(2, 3, 4, 5)
$cols = @("D", "I", "J", "K", "L", "M", "N", "P", "Q")

# Snapshot current (pre-edit) values for the columns that move.
$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Source row for each destination row (rotate up, wrap last to first).
$sourceFor = @{ 2 = 3; 3 = 4; 4 = 5; 5 = 2 }

foreach ($destRow in $rows) {
    $srcRow = $sourceFor[$destRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $snapshot[$srcRow][$c]
    }
}
